$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeByName($slide, $name) {
    $count = $slide.Shapes.Count
    for ($i = 1; $i -le $count; $i++) {
        $candidate = $slide.Shapes.Item($i)
        if ($candidate.Name -eq $name) {
            return $candidate
        }
    }
    return $slide.Shapes.Item($name)
}

# --- Edit 1: "Challenges" textbox (shape "TextBox 5") ---
# Original single run text:
#   "Android Development – Learning and implementing UI best practices in a mobile environment"
# Target (split into 3 runs, the word "UI " removed):
#   Run1: "Android Development – Learning "
#   Run2: "and implementing "
#   Run3: "best practices in a mobile environment"
$shChallenges = Get-ShapeByName $s "TextBox 5"
$trChallenges = $shChallenges.TextFrame.TextRange

# Remove "UI " so "and implementing " is immediately followed by "best practices..."
$uiRange = $trChallenges.Characters(61, 3)
if ($uiRange.Text -eq "UI ") {
    $uiRange.Text = ""
}

# Split "Android Development – Learning " off from "and implementing " by
# re-assigning the same text onto that sub-range (forces a run boundary there).
$firstPart = $trChallenges.Characters(13, 31)
$firstPartExpected = "Android Development " + [char]0x2013 + " Learning "
if ($firstPart.Text -eq $firstPartExpected) {
    $firstPart.Text = $firstPart.Text
}

# --- Edit 2: "Advisor" textbox (shape "TextBox 19") ---
# "Computer Science?" -> "Computer Science" (drop the trailing "?")
$shAdvisor = Get-ShapeByName $s "TextBox 19"
$trAdvisor = $shAdvisor.TextFrame.TextRange
$qLen = $trAdvisor.Length
$qMark = $trAdvisor.Characters($qLen, 1)
if ($qMark.Text -eq "?") {
    $qMark.Text = ""
}
